$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; D = 45092; L = "Primera"; M = 110; N = 10000; O = 11000; P = 10455; Q = "`$/bandeja 18 kilos granel"; R = "Provincia de Curicó"; S = 581; T = 18 },
    @{ Row = 3; D = 44358; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 639; T = 18 },
    @{ Row = 4; D = 44299; L = "Primera"; M = 100; N = 10000; O = 11000; P = 10500; Q = "`$/caja 18 kilos granel"; R = "Región del Maule"; S = 583; T = 18 },
    @{ Row = 5; D = 44299; L = "Segunda"; M = 50; N = 9000; O = 9000; P = 9000; Q = "`$/caja 18 kilos granel"; R = "Región del Maule"; S = 500; T = 18 },
    @{ Row = 6; D = 45079; L = "Primera"; M = 270; N = 11000; O = 12000; P = 11444; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 636; T = 18 },
    @{ Row = 7; D = 44363; L = "Primera"; M = 100; N = 9000; O = 10000; P = 9500; Q = "`$/caja 15 kilos empedrada"; R = "Región de O'Higgins"; S = 633; T = 15 },
    @{ Row = 8; D = 44316; L = "Primera"; M = 100; N = 9000; O = 10000; P = 9500; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 528; T = 18 },
    @{ Row = 9; D = 45126; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 806; T = 18 },
    @{ Row = 10; D = 45013; L = "Primera"; M = 100; N = 9000; O = 10000; P = 9500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 528; T = 18 },
    @{ Row = 11; D = 45027; L = "Primera"; M = 100; N = 9000; O = 10000; P = 9500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 528; T = 18 },
    @{ Row = 12; D = 45037; L = "Primera"; M = 250; N = 9000; O = 9500; P = 9200; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 511; T = 18 },
    @{ Row = 13; D = 45128; L = "Primera"; M = 50; N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 667; T = 18 },
    @{ Row = 14; D = 44272; L = "Primera"; M = 100; N = 9000; O = 10000; P = 9500; Q = "`$/caja 15 kilos granel"; R = "Región de O'Higgins"; S = 633; T = 15 },
    @{ Row = 15; D = 44272; L = "Segunda"; M = 50; N = 8000; O = 8000; P = 8000; Q = "`$/caja 15 kilos granel"; R = "Región de O'Higgins"; S = 533; T = 15 },
    @{ Row = 16; D = 45029; L = "Primera"; M = 100; N = 9000; O = 10000; P = 9500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 528; T = 18 },
    @{ Row = 17; D = 45107; L = "Primera"; M = 100; N = 11000; O = 11000; P = 11000; Q = "`$/caja 18 kilos empedrada"; R = "Región del Maule"; S = 611; T = 18 },
    @{ Row = 18; D = 44698; L = "Primera"; M = 50; N = 10000; O = 10000; P = 10000; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 556; T = 18 },
    @{ Row = 19; D = 45041; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 639; T = 18 },
    @{ Row = 20; D = 45154; L = "Primera"; M = 100; N = 13000; O = 14000; P = 13500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 750; T = 18 },
    @{ Row = 21; D = 45034; L = "Primera"; M = 220; N = 8500; O = 9000; P = 8727; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 485; T = 18 },
    @{ Row = 22; D = 44776; L = "Primera"; M = 50; N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 556; T = 18 },
    @{ Row = 23; D = 44776; L = "Segunda"; M = 50; N = 8000; O = 8000; P = 8000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 444; T = 18 },
    @{ Row = 24; D = 44999; L = "Primera"; M = 100; N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 667; T = 18 },
    @{ Row = 25; D = 44999; L = "Segunda"; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 556; T = 18 },
    @{ Row = 26; D = 44425; L = "Primera"; M = 100; N = 12000; O = 13000; P = 12500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 694; T = 18 },
    @{ Row = 27; D = 44307; L = "Primera"; M = 50; N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 556; T = 18 },
    @{ Row = 28; D = 44307; L = "Segunda"; M = 50; N = 8000; O = 8000; P = 8000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 444; T = 18 },
    @{ Row = 29; D = 45076; L = "Primera"; M = 150; N = 10000; O = 11000; P = 10467; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 582; T = 18 },
    @{ Row = 30; D = 45140; L = "Primera"; M = 50; N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 667; T = 18 },
    @{ Row = 31; D = 45050; L = "Primera"; M = 140; N = 11000; O = 12000; P = 11429; Q = "`$/caja 18 kilos empedrada"; R = "Región de O'Higgins"; S = 635; T = 18 },
    @{ Row = 32; D = 45014; L = "Primera"; M = 100; N = 9000; O = 10000; P = 9500; Q = "`$/bandeja 18 kilos granel"; R = "Región de O'Higgins"; S = 528; T = 18 }
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 12).Value = $item.L
    $ws.Cells.Item($item.Row, 13).Value = $item.M
    $ws.Cells.Item($item.Row, 14).Value = $item.N
    $ws.Cells.Item($item.Row, 15).Value = $item.O
    $ws.Cells.Item($item.Row, 16).Value = $item.P
    $ws.Cells.Item($item.Row, 17).Value = $item.Q
    $ws.Cells.Item($item.Row, 18).Value = $item.R
    $ws.Cells.Item($item.Row, 19).Value = $item.S
    $ws.Cells.Item($item.Row, 20).Value = $item.T
}
